$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 25,14
$arr[0,0] = 12.73908466666667
$arr[0,1] = 38.217254
$arr[0,2] = 0.005953388968763418
$arr[0,3] = 0.006105597140986208
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 162.7225033333333
$arr[0,7] = 488.16751
$arr[0,8] = 0.5231437953541009
$arr[0,9] = 0.5247717033381212
$arr[0,10] = 2072.935747135282
$arr[0,11] = 18656.42172421754
$arr[0,12] = 0.003114478500338132
$arr[0,13] = 0.003204044611571695
$arr[1,0] = 12.73908466666667
$arr[1,1] = 38.217254
$arr[1,2] = 0.005953388968763418
$arr[1,3] = 0.006105597140986208
$arr[1,4] = 1
$arr[1,5] = 0.3333333333333333
$arr[1,6] = 0.2899643333333333
$arr[1,7] = 0.869893
$arr[1,8] = 0.0009322191998643353
$arr[1,9] = 0.0009351200601857102
$arr[1,10] = 3.693880192646889
$arr[1,11] = 33.244921733822
$arr[1,12] = 0.000005549863500941794
$arr[1,13] = 0.000005709466365948723
$arr[2,0] = 12.73908466666667
$arr[2,1] = 38.217254
$arr[2,2] = 0.005953388968763418
$arr[2,3] = 0.006105597140986208
$arr[2,4] = 3
$arr[2,5] = 1
$arr[2,6] = 61.580654
$arr[2,7] = 184.741962
$arr[2,8] = 0.1979783766474813
$arr[2,9] = 0.1985944416431287
$arr[2,10] = 784.4811651347055
$arr[2,11] = 7060.330486212349
$arr[2,12] = 0.001178642283586804
$arr[2,13] = 0.001212537655112039
$arr[3,0] = 12.73908466666667
$arr[3,1] = 38.217254
$arr[3,2] = 0.005953388968763418
$arr[3,3] = 0.006105597140986208
$arr[3,4] = 2
$arr[3,5] = 1
$arr[3,6] = 2.8947245
$arr[3,7] = 5.789449
$arr[3,8] = 0.009306378223129816
$arr[3,9] = 0.00622355841157717
$arr[3,10] = 36.87614049217434
$arr[3,11] = 221.256842953046
$arr[3,12] = 0.00005540448945272114
$arr[3,13] = 0.00003799854044448624
$arr[4,0] = 12.73908466666667
$arr[4,1] = 38.217254
$arr[4,2] = 0.005953388968763418
$arr[4,3] = 0.006105597140986208
$arr[4,4] = 3
$arr[4,5] = 1
$arr[4,6] = 83.559527
$arr[4,7] = 250.678581
$arr[4,8] = 0.2686392305754237
$arr[4,9] = 0.2694751765469873
$arr[4,10] = 1064.471889159619
$arr[4,11] = 9580.247002436576
$arr[4,12] = 0.001599313831884819
$arr[4,13] = 0.001645306867492039
$arr[5,0] = 1964.941406333333
$arr[5,1] = 5894.824219
$arr[5,2] = 0.9182810852447438
$arr[5,3] = 0.9417584502053091
$arr[5,4] = 3
$arr[5,5] = 1
$arr[5,6] = 162.7225033333333
$arr[5,7] = 488.16751
$arr[5,8] = 0.5231437953541009
$arr[5,9] = 0.5247717033381212
$arr[5,10] = 319740.1845418805
$arr[5,11] = 2877661.660876925
$arr[5,12] = 0.480393052136818
$arr[5,13] = 0.4942081860473093
$arr[6,0] = 1964.941406333333
$arr[6,1] = 5894.824219
$arr[6,2] = 0.9182810852447438
$arr[6,3] = 0.9417584502053091
$arr[6,4] = 1
$arr[6,5] = 0.3333333333333333
$arr[6,6] = 0.2899643333333333
$arr[6,7] = 0.869893
$arr[6,8] = 0.0009322191998643353
$arr[6,9] = 0.0009351200601857102
$arr[6,10] = 569.7629249265075
$arr[6,11] = 5127.866324338567
$arr[6,12] = 0.0008560392585374086
$arr[6,13] = 0.0008806572186363898
$arr[7,0] = 1964.941406333333
$arr[7,1] = 5894.824219
$arr[7,2] = 0.9182810852447438
$arr[7,3] = 0.9417584502053091
$arr[7,4] = 3
$arr[7,5] = 1
$arr[7,6] = 61.580654
$arr[7,7] = 184.741962
$arr[7,8] = 0.1979783766474813
$arr[7,9] = 0.1985944416431287
$arr[7,10] = 121002.3768736864
$arr[7,11] = 1089021.391863178
$arr[7,12] = 0.1817997985628418
$arr[7,13] = 0.1870279935812216
$arr[8,0] = 1964.941406333333
$arr[8,1] = 5894.824219
$arr[8,2] = 0.9182810852447438
$arr[8,3] = 0.9417584502053091
$arr[8,4] = 2
$arr[8,5] = 1
$arr[8,6] = 2.8947245
$arr[8,7] = 5.789449
$arr[8,8] = 0.009306378223129816
$arr[8,9] = 0.00622355841157717
$arr[8,10] = 5687.964029977556
$arr[8,11] = 34127.78417986533
$arr[8,12] = 0.008545871094433697
$arr[8,13] = 0.005861088724449131
$arr[9,0] = 1964.941406333333
$arr[9,1] = 5894.824219
$arr[9,2] = 0.9182810852447438
$arr[9,3] = 0.9417584502053091
$arr[9,4] = 3
$arr[9,5] = 1
$arr[9,6] = 83.559527
$arr[9,7] = 250.678581
$arr[9,8] = 0.2686392305754237
$arr[9,9] = 0.2694751765469873
$arr[9,10] = 164189.5744959282
$arr[9,11] = 1477706.170463353
$arr[9,12] = 0.246686324192113
$arr[9,13] = 0.2537805246336928
$arr[10,0] = 1.091866333333334
$arr[10,1] = 3.275599000000001
$arr[10,2] = 0.0005102646818291153
$arr[10,3] = 0.0005233104369407934
$arr[10,4] = 3
$arr[10,5] = 1
$arr[10,6] = 162.7225033333333
$arr[10,7] = 488.16751
$arr[10,8] = 0.5231437953541009
$arr[10,9] = 0.5247717033381212
$arr[10,10] = 177.6712230653878
$arr[10,11] = 1599.04100758849
$arr[10,12] = 0.0002669418022872361
$arr[10,13] = 0.0002746185093680366
$arr[11,0] = 1.091866333333334
$arr[11,1] = 3.275599000000001
$arr[11,2] = 0.0005102646818291153
$arr[11,3] = 0.0005233104369407934
$arr[11,4] = 1
$arr[11,5] = 0.3333333333333333
$arr[11,6] = 0.2899643333333333
$arr[11,7] = 0.869893
$arr[11,8] = 0.0009322191998643353
$arr[11,9] = 0.0009351200601857102
$arr[11,10] = 0.3166022934341112
$arr[11,11] = 2.849420640907001
$arr[11,12] = 0.0000004756785334137675
$arr[11,13] = 0.0000004893580872878849
$arr[12,0] = 1.091866333333334
$arr[12,1] = 3.275599000000001
$arr[12,2] = 0.0005102646818291153
$arr[12,3] = 0.0005233104369407934
$arr[12,4] = 3
$arr[12,5] = 1
$arr[12,6] = 61.580654
$arr[12,7] = 184.741962
$arr[12,8] = 0.1979783766474813
$arr[12,9] = 0.1985944416431287
$arr[12,10] = 67.23784288724869
$arr[12,11] = 605.1405859852381
$arr[12,12] = 0.0001010213733690718
$arr[12,13] = 0.0001039265440302786
$arr[13,0] = 1.091866333333334
$arr[13,1] = 3.275599000000001
$arr[13,2] = 0.0005102646818291153
$arr[13,3] = 0.0005233104369407934
$arr[13,4] = 2
$arr[13,5] = 1
$arr[13,6] = 2.8947245
$arr[13,7] = 5.789449
$arr[13,8] = 0.009306378223129816
$arr[13,9] = 0.00622355841157717
$arr[13,10] = 3.160652225825168
$arr[13,11] = 18.963913354951
$arr[13,12] = 0.000004748716123006742
$arr[13,13] = 0.000003256853071688999
$arr[14,0] = 1.091866333333334
$arr[14,1] = 3.275599000000001
$arr[14,2] = 0.0005102646818291153
$arr[14,3] = 0.0005233104369407934
$arr[14,4] = 3
$arr[14,5] = 1
$arr[14,6] = 83.559527
$arr[14,7] = 250.678581
$arr[14,8] = 0.2686392305754237
$arr[14,9] = 0.2694751765469873
$arr[14,10] = 91.23583436055769
$arr[14,11] = 821.1225092450192
$arr[14,12] = 0.0001370771115163869
$arr[14,13] = 0.0001410191723835013
$arr[15,0] = 160.0313415
$arr[15,1] = 320.062683
$arr[15,2] = 0.0747878554913321
$arr[15,3] = 0.05113328661083746
$arr[15,4] = 3
$arr[15,5] = 1
$arr[15,6] = 162.7225033333333
$arr[15,7] = 488.16751
$arr[15,8] = 0.5231437953541009
$arr[15,9] = 0.5247717033381212
$arr[15,10] = 26040.70050067155
$arr[15,11] = 156244.2030040293
$arr[15,12] = 0.03912480256812952
$arr[15,13] = 0.02683330191204552
$arr[16,0] = 160.0313415
$arr[16,1] = 320.062683
$arr[16,2] = 0.0747878554913321
$arr[16,3] = 0.05113328661083746
$arr[16,4] = 1
$arr[16,5] = 0.3333333333333333
$arr[16,6] = 0.2899643333333333
$arr[16,7] = 0.869893
$arr[16,8] = 0.0009322191998643353
$arr[16,9] = 0.0009351200601857102
$arr[16,10] = 46.4033812504865
$arr[16,11] = 278.420287502919
$arr[16,12] = 0.00006971867480569915
$arr[16,13] = 0.0000478157620530195
$arr[17,0] = 160.0313415
$arr[17,1] = 320.062683
$arr[17,2] = 0.0747878554913321
$arr[17,3] = 0.05113328661083746
$arr[17,4] = 3
$arr[17,5] = 1
$arr[17,6] = 61.580654
$arr[17,7] = 184.741962
$arr[17,8] = 0.1979783766474813
$arr[17,9] = 0.1985944416431287
$arr[17,10] = 9854.834670067341
$arr[17,11] = 59129.00802040404
$arr[17,12] = 0.01480637822312035
$arr[17,13] = 0.01015478650385734
$arr[18,0] = 160.0313415
$arr[18,1] = 320.062683
$arr[18,2] = 0.0747878554913321
$arr[18,3] = 0.05113328661083746
$arr[18,4] = 2
$arr[18,5] = 1
$arr[18,6] = 2.8947245
$arr[18,7] = 5.789449
$arr[18,8] = 0.009306378223129816
$arr[18,9] = 0.00622355841157717
$arr[18,10] = 463.2466450079168
$arr[18,11] = 1852.986580031667
$arr[18,12] = 0.0006960040696991127
$arr[18,13] = 0.0003182309959984638
$arr[19,0] = 160.0313415
$arr[19,1] = 320.062683
$arr[19,2] = 0.0747878554913321
$arr[19,3] = 0.05113328661083746
$arr[19,4] = 3
$arr[19,5] = 1
$arr[19,6] = 83.559527
$arr[19,7] = 250.678581
$arr[19,8] = 0.2686392305754237
$arr[19,9] = 0.2694751765469873
$arr[19,10] = 13372.14320091547
$arr[19,11] = 80232.85920549283
$arr[19,12] = 0.02009095195557743
$arr[19,13] = 0.01377915143688313
$arr[20,0] = 1.000156333333333
$arr[20,1] = 3.000469
$arr[20,2] = 0.0004674056133315229
$arr[20,3] = 0.0004793556059265206
$arr[20,4] = 3
$arr[20,5] = 1
$arr[20,6] = 162.7225033333333
$arr[20,7] = 488.16751
$arr[20,8] = 0.5231437953541009
$arr[20,9] = 0.5247717033381212
$arr[20,10] = 162.7479422846878
$arr[20,11] = 1464.73148056219
$arr[20,12] = 0.0002445203465280643
$arr[20,13] = 0.0002515522578267374
$arr[21,0] = 1.000156333333333
$arr[21,1] = 3.000469
$arr[21,2] = 0.0004674056133315229
$arr[21,3] = 0.0004793556059265206
$arr[21,4] = 1
$arr[21,5] = 0.3333333333333333
$arr[21,6] = 0.2899643333333333
$arr[21,7] = 0.869893
$arr[21,8] = 0.0009322191998643353
$arr[21,9] = 0.0009351200601857102
$arr[21,10] = 0.2900096644241111
$arr[21,11] = 2.610086979817
$arr[21,12] = 0.0000004357244868720112
$arr[21,13] = 0.0000004482550430643655
$arr[22,0] = 1.000156333333333
$arr[22,1] = 3.000469
$arr[22,2] = 0.0004674056133315229
$arr[22,3] = 0.0004793556059265206
$arr[22,4] = 3
$arr[22,5] = 1
$arr[22,6] = 61.580654
$arr[22,7] = 184.741962
$arr[22,8] = 0.1979783766474813
$arr[22,9] = 0.1985944416431287
$arr[22,10] = 61.59028110890866
$arr[22,11] = 554.3125299801779
$arr[22,12] = 0.00009253620456329525
$arr[22,13] = 0.000095197358907481
$arr[23,0] = 1.000156333333333
$arr[23,1] = 3.000469
$arr[23,2] = 0.0004674056133315229
$arr[23,3] = 0.0004793556059265206
$arr[23,4] = 2
$arr[23,5] = 1
$arr[23,6] = 2.8947245
$arr[23,7] = 5.789449
$arr[23,8] = 0.009306378223129816
$arr[23,9] = 0.00622355841157717
$arr[23,10] = 2.895177041930166
$arr[23,11] = 17.371062251581
$arr[23,12] = 0.00000434985342127712
$arr[23,13] = 0.000002983297613400669
$arr[24,0] = 1.000156333333333
$arr[24,1] = 3.000469
$arr[24,2] = 0.0004674056133315229
$arr[24,3] = 0.0004793556059265206
$arr[24,4] = 3
$arr[24,5] = 1
$arr[24,6] = 83.559527
$arr[24,7] = 250.678581
$arr[24,8] = 0.2686392305754237
$arr[24,9] = 0.2694751765469873
$arr[24,10] = 83.57259013938766
$arr[24,11] = 752.1533112544889
$arr[24,12] = 0.0001255634843320143
$arr[24,13] = 0.0001291744365358372

$ws.Range("G2:T26").Value = $arr
